$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# --- Title -----------------------------------------------------------
# The H1 heading (paragraph 1) and the bold "title" paragraph near the end
# (paragraph 43) both hold the exact same original string, so a single
# document-wide replace-all updates both occurrences at once while each
# run keeps its own formatting (e.g. the bold rPr on the later copy).
$d.Content.Find.Execute(
    "Play Marvelous Furlongs for Free - Review", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Play Marvelous Furlongs Free - Exciting Slot Game Review", 2) | Out-Null

# --- "What we like" bullets (paragraphs 36-39) ------------------------
# Each bullet keeps its paragraph slot, but the wording shifts between
# slots, and some new values equal other bullets' old values (e.g. bullet
# 38's new text equals bullet 37's old text) - so every replace is scoped
# to its own paragraph's Range to avoid cross-matching, and done in an
# order where a not-yet-edited "old" string is never reintroduced earlier.
Replace-InParagraph 36 "Stacked Respins feature" "High volatility and excellent RTP"
Replace-InParagraph 37 "Two different Free Spin features" "Stacked Respin feature for more chances to win"
Replace-InParagraph 38 "Excellent RTP of 96.83%" "Two different Free Spin features"
Replace-InParagraph 39 "High-quality graphics and design" "Stunning design and graphics"

# --- "What we don't like" bullets (paragraph 42) ----------------------
Replace-InParagraph 42 "High volatility might not be for everyone" "Limited betting options"

# --- Meta description (italic paragraph at the end) -------------------
$d.Content.Find.Execute(
    "Read our review of Marvelous Furlongs, a high-volatility online slot game with Stacked Respins and two Free Spins features. Play the game for free.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Marvelous Furlongs, an exciting slot game with high volatility and two Free Spin features that you can play for free.", 2) | Out-Null

Write-Output "Edit complete"
